$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 2 (shifts existing data rows 2..22 down to 3..23)
$ws.Rows.Item(2).Insert()

# New row 2 should contain only C2 = "---"; the insert drags E1's format into E2, clear it
$ws.Cells.Item(2, 3).Value = "---"
$ws.Cells.Item(2, 5).Clear()

# Apply number format "@" (text, numFmtId 49) across column C (header + all data rows)
$ws.Range("C1:C23").NumberFormat = "@"

# Selection moves to C3
$ws.Range("C3").Select()
